$wb = $excel.ActiveWorkbook

# Both the "展览" sheet and the "全部类型" sheet carry the same event data
# and need the same numeric updates (F2: 153 -> 155, F3: 29 -> 30).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 155
    $ws.Range("F3").Value = 30
}
